# Weekly crime-stat refresh: bump the report volume/number and the
# covered week dates, and replace the precinct figures for the week
# (and the rows whose 28-day / YTD / 2-yr roll-ups depend on them)
# with the newly collected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Masthead: "Volume 30   Number  44" -> "...  45" --------------
$ws.Range("A8").Value = "Volume 30   Number  45"

# ---- "Report Covering the Week  10/30/2023  Through  11/5/2023" ---
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# ---- Row 15 (Rape) -------------------------------------------------
# F15 goes from a real count (2) to the "no data" marker used by C15/D15.
$ws.Range("C15").Copy($ws.Range("F15"))
$ws.Range("H15").Value = -100
$ws.Range("M15").Value = 57.142857142857
$ws.Range("N15").Value = 100

# ---- Row 16 (Robbery) ----------------------------------------------
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 23
$ws.Range("H16").Value = 109.090909090909
$ws.Range("I16").Value = 230
$ws.Range("J16").Value = 169
$ws.Range("K16").Value = 36.094674556213
$ws.Range("L16").Value = 123.300970873786
$ws.Range("M16").Value = 36.904761904761
$ws.Range("N16").Value = -73.922902494331

# ---- Row 17 (Fel. Assault) ------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = -5.882352941176
$ws.Range("I17").Value = 214
$ws.Range("J17").Value = 195
$ws.Range("K17").Value = 9.743589743589
$ws.Range("L17").Value = 18.232044198895
$ws.Range("M17").Value = 91.071428571428
$ws.Range("N17").Value = -21.323529411764

# ---- Row 18 (Burglary) ----------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 53.846153846153
$ws.Range("I18").Value = 186
$ws.Range("J18").Value = 147
$ws.Range("K18").Value = 26.530612244898
$ws.Range("L18").Value = 14.814814814814
$ws.Range("M18").Value = -13.084112149532
$ws.Range("N18").Value = -85.993975903614

# ---- Row 19 (Gr. Larceny) --------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 6.521739130434
$ws.Range("I19").Value = 642
$ws.Range("J19").Value = 608
$ws.Range("K19").Value = 5.592105263157
$ws.Range("L19").Value = 65.463917525773
$ws.Range("M19").Value = 58.128078817734
$ws.Range("N19").Value = -19.346733668341

# ---- Row 20 (G.L.A.) --------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -28.571428571428
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 8.333333333333
$ws.Range("I20").Value = 271
$ws.Range("J20").Value = 208
$ws.Range("K20").Value = 30.288461538461
$ws.Range("L20").Value = 69.375
$ws.Range("M20").Value = 49.723756906077
$ws.Range("N20").Value = -85.191256830601

# ---- Row 21 (TOTAL) ----------------------------------------------------
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -9.677419354838
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = 19.642857142857
$ws.Range("I21").Value = 1565
$ws.Range("J21").Value = 1346
$ws.Range("K21").Value = 16.270430906389
$ws.Range("L21").Value = 55.412115193644
$ws.Range("M21").Value = 42.661804922516
$ws.Range("N21").Value = -69.475326701774

# ---- Row 22 (Transit) ---------------------------------------------------
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -11.111111111111
$ws.Range("I22").Value = 76
$ws.Range("J22").Value = 72
$ws.Range("K22").Value = 5.555555555555
$ws.Range("L22").Value = 123.529411764706
$ws.Range("M22").Value = 72.727272727272

# ---- Row 24 (Petit Larceny) ---------------------------------------------
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = 6.060606060606
$ws.Range("F24").Value = 198
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = 47.761194029850
$ws.Range("I24").Value = 1786
$ws.Range("J24").Value = 1307
$ws.Range("K24").Value = 36.648814078041
$ws.Range("L24").Value = 61.336946702800
$ws.Range("M24").Value = 117.010935601458

# ---- Row 25 (Misd. Assault) ----------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -35.714285714285
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -29.090909090909
$ws.Range("I25").Value = 450
$ws.Range("J25").Value = 492
$ws.Range("K25").Value = -8.536585365853
$ws.Range("L25").Value = 8.695652173913
$ws.Range("M25").Value = 3.211009174311

# ---- Row 26 (UCR Rape*) ----------------------------------------------------
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50
$ws.Range("L26").Value = 45

# ---- Row 27 (Other Sex Crimes) ---------------------------------------------
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = -11.111111111111
$ws.Range("I27").Value = 76
$ws.Range("J27").Value = 80
$ws.Range("K27").Value = -5
$ws.Range("L27").Value = 43.396226415094

# ---- Row 30 (Hate Crimes) ----------------------------------------------------
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 8
$ws.Range("K30").Value = -11.111111111111
$ws.Range("L30").Value = -33.333333333333
